# Remove the AutoFilter criterion on column E ("name") and re-sort the
# data range by column A ascending, as part of tidying up the cleaning
# documentation sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trips_p1_stations")

# Clear any active filter criteria (this also unhides any rows that were
# hidden by the filter).
if ($ws.FilterMode) {
    $ws.ShowAllData()
}

# Re-apply sort: sort A2:E150 by column A ascending.
$sortRange = $ws.Range("A1:E150")
$key1 = $ws.Range("A1:A150")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key1, 0, 1, 0, 0) | Out-Null

$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Update the active selection as in the final saved file.
$ws.Range("C7").Select()
